$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "u_entry" header in column G, row 1 (matches header style of other columns, no explicit style)
$ws.Range("G1").Value = "u_entry"

# Update the selected/active cell to G2 (as seen in the target sheetView selection)
$ws.Range("G2").Select()
